{"js": "// Charms & Treasures review \u2014 rewrite title, meta description, and the\n// \"What we like\" / \"What we don't like\" bullet lists.\n//\n// Each old string below is unique within the document, so a scoped\n// search-and-replace (matchCase, no wildcards) safely targets exactly the\n// run that needs to change without disturbing surrounding formatting.\n//\n// NOTE on ordering: \"High-paying symbols for significant payouts\" is both\n// the OLD text of one bullet and the NEW text of another bullet, so that\n// pair is processed first to avoid a renamed paragraph being matched again\n// by a later replacement.\n\nconst body = context.document.body;\n\nconst replacements = [\n  // \"What we like\" bullet 3 -> must run before the bullet-1 replacement\n  // below (whose new text equals this rule's old text).\n  [\n    \"High-paying symbols for significant payouts\",\n    \"Autospin feature for convenient gameplay\",\n  ],\n  // Title (appears in the Heading1 and again in the bold \"title\" paragraph\n  // near the end of the document) \u2014 replace every occurrence.\n  [\n    \"Play Charms & Treasures Slot Free | Review of 5-Reel Game\",\n    \"Play Charms & Treasures Free: Review of Exciting Slot Game\",\n  ],\n  // \"What we like\" bullets\n  [\n    \"Bonus game exponentially increases possibility of high-paying combos\",\n    \"High-paying symbols for significant payouts\",\n  ],\n  [\n    \"Autospin feature allows for defined amount of automatic spins\",\n    \"Bonus game with a special grid for high-paying combinations\",\n  ],\n  [\n    \"Affordable minimum bet of \\u20ac0.25 for all types of players\",\n    \"Well-executed theme and visuals\",\n  ],\n  // \"What we don't like\" bullets\n  [\n    \"Limited availability of free spins, only triggered by Scatter symbol\",\n    \"High volatility may not appeal to all players\",\n  ],\n  [\n    \"High volatility may lead to frequent losses\",\n    \"Limited number of Free Spins from the Scatter symbol\",\n  ],\n  // Meta description (italic paragraph at the very end)\n  [\n    \"Read our review of Charms & Treasures slot, a 5-reel game with an Autospin feature and Bonus game. Play for free and enjoy significant payouts!\",\n    \"Discover the gameplay features and bonus potential of Charms & Treasures. Play for free now!\",\n  ],\n];\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  for (const r of results.items) {\n    r.insertText(newText, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "# Charms & Treasures review \u2014 rewrite title, meta description, and the\n# \"What we like\" / \"What we don't like\" bullet lists.\n#\n# Each old string is unique within the document, so a scoped Find/Replace\n# (MatchCase, no wildcards) safely targets exactly the run that needs to\n# change without disturbing surrounding formatting.\n#\n# NOTE on ordering: \"High-paying symbols for significant payouts\" is both\n# the OLD text of one bullet and the NEW text of another bullet, so that\n# pair is processed first to avoid a renamed paragraph being matched again\n# by a later replacement.\n\n$d = $word.ActiveDocument\n\nfunction Replace-Text($old, $new) {\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $old\n    $find.Replacement.Text = $new\n    # wdFindContinue = 1, wdReplaceAll = 2\n    $find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)\n}\n\n# \"What we like\" bullet 3 -> must run before the bullet-1 replacement below\n# (whose new text equals this rule's old text).\nReplace-Text \"High-paying symbols for significant payouts\" \"Autospin feature for convenient gameplay\"\n\n# Title (appears in the Heading1 and again in the bold \"title\" paragraph\n# near the end of the document) \u2014 replace every occurrence.\nReplace-Text \"Play Charms & Treasures Slot Free | Review of 5-Reel Game\" \"Play Charms & Treasures Free: Review of Exciting Slot Game\"\n\n# \"What we like\" bullets\nReplace-Text \"Bonus game exponentially increases possibility of high-paying combos\" \"High-paying symbols for significant payouts\"\nReplace-Text \"Autospin feature allows for defined amount of automatic spins\" \"Bonus game with a special grid for high-paying combinations\"\nReplace-Text \"Affordable minimum bet of \u20ac0.25 for all types of players\" \"Well-executed theme and visuals\"\n\n# \"What we don't like\" bullets\nReplace-Text \"Limited availability of free spins, only triggered by Scatter symbol\" \"High volatility may not appeal to all players\"\nReplace-Text \"High volatility may lead to frequent losses\" \"Limited number of Free Spins from the Scatter symbol\"\n\n# Meta description (italic paragraph at the very end)\nReplace-Text \"Read our review of Charms & Treasures slot, a 5-reel game with an Autospin feature and Bonus game. Play for free and enjoy significant payouts!\" \"Discover the gameplay features and bonus potential of Charms & Treasures. Play for free now!\"\n"}
